$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old trailing columns E:F that are no longer used
$ws.Range("E1:F2").ClearContents()

# Header row (bold), now only 4 columns: Nama Peserta, Tanggal, Listening, Reading
$ws.Range("A1").Value = "Nama Peserta"
$ws.Range("B1").Value = "Tanggal"
$ws.Range("C1").Value = "Listening"
$ws.Range("D1").Value = "Reading"
$ws.Range("A1:D1").Font.Bold = $true

# Pre-apply the date number format to B2:B6 BEFORE writing values so the
# engine doesn't mint a transient/orphan auto-date style first.
$ws.Range("B2:B6").NumberFormat = "yyyy-mm-dd"

# Data rows
$ws.Range("A2").Value = "Lyra Faiqah Bilqis"
$ws.Range("B2").Value = "2025-04-15"
$ws.Range("C2").Value = 250
$ws.Range("D2").Value = 260

$ws.Range("A3").Value = "Satria Rakhmadani"
$ws.Range("B3").Value = "2025-04-13"
$ws.Range("C3").Value = 200
$ws.Range("D3").Value = 180

$ws.Range("A4").Value = "Dimas Wahyu"
$ws.Range("B4").Value = "2025-04-15"
$ws.Range("C4").Value = 300
$ws.Range("D4").Value = 120

$ws.Range("A5").Value = "Nicholas Saputra"
$ws.Range("B5").Value = "2025-04-15"
$ws.Range("C5").Value = 300
$ws.Range("D5").Value = 210

$ws.Range("A6").Value = "Mamat Alkatiri"
$ws.Range("B6").Value = "2025-04-13"
$ws.Range("C6").Value = 200
$ws.Range("D6").Value = 250

$ws.Columns.Item(1).ColumnWidth = 26.8571428571429
$ws.Columns.Item(2).ColumnWidth = 17.2857142857143
$ws.Columns.Item(6).ColumnWidth = 13.8571428571429

$ws.Range("C7").Select()
